$d = $word.ActiveDocument

function Find-ParagraphIndex($doc, $text) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $p = $doc.Paragraphs.Item($i)
        if ($p.Range.Text.TrimEnd() -eq $text) {
            return $i
        }
    }
    return -1
}

# --- Edit 1: remove the whole paragraph about internal notifications ---
$notifText = "El sistema debe permitir enviar notificaciones internas al docente cuando un estudiante cambie de nivel de riesgo."
$notifIdx = Find-ParagraphIndex $d $notifText
if ($notifIdx -gt 0) {
    $d.Paragraphs.Item($notifIdx).Range.Delete()
}

# --- Edit 2: blank out the "reportes globales" paragraph and add a new blank
#             paragraph right after it (both without bullet numbering, using
#             the same hanging-indent formatting) ---
$reportsText = "El sistema debe permitir que los administradores generen reportes globales de toda la institución."
$reportsIdx = Find-ParagraphIndex $d $reportsText
if ($reportsIdx -gt 0) {
    $p = $d.Paragraphs.Item($reportsIdx)
    $rng = $p.Range
    [void]$rng.MoveEnd(1, -1)
    $found = $rng.Find.Execute($reportsText, $true, $false, $false, $false, $false, $true, 1, $false, "^p", 2)

    $p1 = $d.Paragraphs.Item($reportsIdx)
    $p2 = $d.Paragraphs.Item($reportsIdx + 1)

    foreach ($pp in @($p1, $p2)) {
        $pp.Range.ListFormat.RemoveNumbers()
        $pp.Range.ParagraphFormat.LeftIndent = 18
        $pp.Range.ParagraphFormat.FirstLineIndent = -18
    }
}
